$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column H (Generated Date (UTC)) / column I (Generated By)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("H$r").Value = "2025-05-01 10:03:02"
    $ws.Range("I$r").Value = "wolketichif"
}
